$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 99.2
$ws.Range("I2").Value = 99.2
$ws.Range("K2").Value = 99.2
$ws.Range("M2").Value = 13.8

$ws.Range("H9").Value = 62
$ws.Range("J9").Value = 25
$ws.Range("L9").Value = 25
$ws.Range("N9").Value = -363

$ws.Range("H113").Value = 2077.6
$ws.Range("I113").Value = 1397.5
$ws.Range("J113").Value = 2531
$ws.Range("K113").Value = 1397.5
$ws.Range("L113").Value = 2531
$ws.Range("M113").Value = 1856.5
$ws.Range("N113").Value = -9039

$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").Value = $null

$ws.Range("H138").Value = 0
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("M138").Value = $null
$ws.Range("N138").Value = $null

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("M30").Value = $null

$ws.Range("H61").Value = 2933
$ws.Range("I61").Value = 2900
$ws.Range("J61").Value = 2999
$ws.Range("K61").Value = 2900
$ws.Range("L61").Value = 2999
$ws.Range("M61").Value = -2688
$ws.Range("N61").Value = -3423

$ws.Range("H74").Value = 3971.75
$ws.Range("I74").Value = 3971.75
$ws.Range("K74").Value = 3971.75
$ws.Range("M74").Value = -3097.75

$ws.Range("H77").Value = 3971.75
$ws.Range("I77").Value = 3971.75
$ws.Range("K77").Value = 19858.75
$ws.Range("M77").Value = -15490.75

$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").Value = $null

$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = $null
$ws.Range("N132").Value = $null

$ws.Range("H136").Value = 2933
$ws.Range("I136").Value = 2900
$ws.Range("J136").Value = 2999
$ws.Range("K136").Value = 8700
$ws.Range("L136").Value = 8997
$ws.Range("M136").Value = -6150
$ws.Range("N136").Value = -14097

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = $null
$ws.Range("N8").Value = $null

$ws.Range("H22").Value = 298.2
$ws.Range("I22").Value = 298.22223
$ws.Range("K22").Value = 298.22223
$ws.Range("M22").Value = 51.77776999999998

$ws.Range("H31").Value = 3500
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 3500
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 3500
$ws.Range("M31").Value = $null
$ws.Range("N31").Value = -4090

$ws.Range("H34").Value = 3500
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 3500
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 3500
$ws.Range("M34").Value = $null
$ws.Range("N34").Value = -3904

$ws.Range("H58").Value = 3909.5833
$ws.Range("I58").Value = 3224.2222
$ws.Range("K58").Value = 3224.2222
$ws.Range("M58").Value = -3021.2222

$ws.Range("H104").Value = 17285
$ws.Range("J104").Value = 17285
$ws.Range("L104").Value = 17285
$ws.Range("N104").Value = -22527

$ws.Range("H134").Value = 5632.3335
$ws.Range("I134").Value = 3430.5
$ws.Range("K134").Value = 10291.5
$ws.Range("M134").Value = -7756.5

$ws.Range("H136").Value = 3909.5833
$ws.Range("I136").Value = 3224.2222
$ws.Range("K136").Value = 9672.6666
$ws.Range("M136").Value = -7122.6666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 827.6667
$ws.Range("J129").Value = 917
$ws.Range("L129").Value = 2751
$ws.Range("N129").Value = -12751

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 26666.666
$ws.Range("J20").Value = 26666.666
$ws.Range("L20").Value = 26666.666
$ws.Range("N20").Value = -27156.666

$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").Value = $null

$ws.Range("H97").Value = 4284.5
$ws.Range("I97").Value = 4284.5
$ws.Range("K97").Value = 4284.5
$ws.Range("M97").Value = -3788.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 639.4
$ws.Range("I22").Value = 599
$ws.Range("K22").Value = 599
$ws.Range("M22").Value = -304

$ws.Range("H25").Value = 307
$ws.Range("I25").Value = 307
$ws.Range("K25").Value = 307
$ws.Range("M25").Value = -77

$ws.Range("H27").Value = 639.4
$ws.Range("I27").Value = 599
$ws.Range("K27").Value = 599
$ws.Range("M27").Value = -492

$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = $null
$ws.Range("N31").Value = $null

$ws.Range("H38").Value = 30000
$ws.Range("J38").Value = 30000
$ws.Range("L38").Value = 30000
$ws.Range("N38").Value = -30820

$ws.Range("H48").Value = 0
$ws.Range("I48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("M48").Value = $null

$ws.Range("H100").Value = 3334.1667
$ws.Range("I100").Value = 2003
$ws.Range("J100").Value = 3600.4
$ws.Range("K100").Value = 2003
$ws.Range("L100").Value = 3600.4
$ws.Range("M100").Value = -1462
$ws.Range("N100").Value = -4682.4

$ws.Range("H106").Value = 8580
$ws.Range("J106").Value = 8580
$ws.Range("L106").Value = 8580
$ws.Range("N106").Value = -11104

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 20000
$ws.Range("I8").Value = 20000
$ws.Range("K8").Value = 20000
$ws.Range("M8").Value = -19860

$ws.Range("H17").Value = 500
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").Value = $null

$ws.Range("H19").Value = 30000000
$ws.Range("J19").Value = 30000000
$ws.Range("L19").Value = 30000000
$ws.Range("N19").Value = -30000348

$ws.Range("H86").Value = 30000
$ws.Range("J86").Value = 30000
$ws.Range("L86").Value = 30000
$ws.Range("N86").Value = -32246

$ws.Range("H89").Value = 30000
$ws.Range("J89").Value = 30000
$ws.Range("L89").Value = 150000
$ws.Range("N89").Value = -161232

$ws.Range("H100").Value = 3037
$ws.Range("I100").Value = 2479.6
$ws.Range("J100").Value = 3966
$ws.Range("K100").Value = 4959.2
$ws.Range("L100").Value = 7932
$ws.Range("M100").Value = -4418.2
$ws.Range("N100").Value = -9014
